$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, pushing existing rows 17:124 down to 18:125.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new data record (dd/mm stored as Excel serial date).
$ws.Cells.Item(17, 1).Value = 5
$ws.Cells.Item(17, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(17, 3).Value = "Maule"
$ws.Cells.Item(17, 4).Value = 44670
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 7
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100108
$ws.Cells.Item(17, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(17, 9).Value = 100108002
$ws.Cells.Item(17, 10).Value = "Mango"
$ws.Cells.Item(17, 11).Value = "Sin especificar"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 340
$ws.Cells.Item(17, 14).Value = 7000
$ws.Cells.Item(17, 15).Value = 7500
$ws.Cells.Item(17, 16).Value = 7294
$ws.Cells.Item(17, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(17, 18).Value = "Ecuador"
$ws.Cells.Item(17, 19).Value = 1824
$ws.Cells.Item(17, 20).Value = 4
